# Apply the PS-7.docx revision: fix typos/spacing/quotes, move sentences
# between paragraphs, and split the final paragraph into two.
$d = $word.ActiveDocument

# The source has 3 paragraphs but the target has 4: the trailing
# "Dated this day of ..." sentence becomes its own paragraph, split
# off of paragraph 3. Insert a paragraph break right before the final
# paragraph mark of paragraph 3 so we end up with an empty 4th paragraph
# to populate.
$p3 = $d.Paragraphs(3)
$splitPoint = $p3.Range.End - 1
$d.Range($splitPoint, $splitPoint).InsertParagraphAfter()

# Rewrite each paragraph's text in full (this also relocates the
# "Under this Guarantee the Vendor is providing..." sentence from the
# end of paragraph 1 to the start of paragraph 2, and the "The Guarantor
# hereby guarantee..." sentence from the end of paragraph 2 to the start
# of paragraph 3, per the revision).
$d.Paragraphs(1).Range.Text = 'BG NO:0005NDLG01739650 Issuance Date: 30-11-2022 BG NO: 0005NDLG01739650 Issuance Date:30-11-2022 BANK GUARANTEE FOR PERFORMANCE GUARANTEE To , Raheja Private Limited , 8 Shroff Chambers , Girgaon , M .G Road Mumbai-400092 This Guarantee is made by Impactsure Bank Limited , having it’s registered office at A207 , Eastern Business District , Bhandup W , Mumbai 400078 (hereinafter referred to as the "GUARANTOR” OR “BANK" which expression shall unless repugnant to the context or meaning thereof include its successors and permitted assigns) in favour of Raheja Private Limited , a public limited company , incorporated and having its registered office at 8 Shroff Chambers , Girgaon , M .G Road , Mumbai-400092 (hereinafter referred to as “Beneficiary”). Whereas M/s .Tango Textiles , having its registered office at Hiranandani , Lake Road , Powai , Mumabi -400615 (hereinafter referred to as "Applicant" which expression shall unless repugnant to the context or meaning thereof include its successors and permitted assigns) having executed an Amendment-1 No .(hereinafter called the Agreement) for Supply , Installation and Maintenance . And Whereas under the terms of the said Agreement and in accordance with the other conditions of Agreement , "Applicant" is required to provide with an irrevocable Bank Guarantee of value equivalent to 5% to Total contract value .'
$d.Paragraphs(2).Range.Text = 'Under this Guarantee the Vendor is providing Performance Bank Guarantee for 225 retail outlets . Based on the above facts “Applicant” has to provide a Bank Guarantee of INR 22 ,000/-(Rupees Twenty Two Thousand Only) Now Therefore this Guarantee witnesses as follows: 1 . In consideration of the aforesaid premises and in consideration of the faithful performance by "Applicant" the terms and conditions of the said Agreement has to be guaranteed by the Bank , we , the Guarantor , hereby guarantee that "Applicant" will duly comply and faithfully perform all his obligation and his responsibility under the said Agreement failing which we , the Guarantor , do undertake to pay to 22 ,000/-(Rupees Twenty Two Thousand Only )immediately on a written demand being made on us , without demur , dispute or objection of whatsoever nature and without reference to "Applicant" such amount or amounts as the Guarantor may be called upon to pay not exceeding in the aggregate a sum of INR 22 ,000/-(Rupees Twenty Two Thousand Only ) Notwithstanding anything to the contrary , Impactsure Bank Limited decision as to whether "Applicant" has made any fault or defaults and the amount to which is BG NO:0005NDLG01739650 Issuance Date: 30-11-2022 entitled by the reason thereof will be binding on us and we shall not be entitled to ask to establish their claims under this Guarantee but , we , the Guarantor will pay the sum forthwith without any objection .'
$d.Paragraphs(3).Range.Text = 'The Guarantor hereby guarantee to the due compliance and observance by Applicant" of the Terms and Conditions of the Agreement and the Guarantor hereby undertakes , that this Guarantee shall be kept valid and binding on the Guarantor from the date of this Guarantee till 30-11-2022 (“Expiry Date”) with a claim period of six months from the Expiry Date i .e . 30-05-2023 and shall not be terminable by notice or any change in the constitution of the Guarantor or by any other reasons whatsoever and the liability of the Guarantor hereunder shall not be impaired or discharged by any extension of the time or variations or alterations made , given , conceded or agreed with or without the Guarantor''s knowledge or consent by or between the parties to the said Agreement . Notwithstanding anything contained herein above: (i) Our liability under this Bank Guarantee shall not exceed INR 22 ,000/-(Rupees Twenty Two Thousand Only ) (ii) This Bank Guarantee shall be valid up to 30-11-2022 post which it ceases to be in effect in all respects whether or not the original Bank Guarantee is returned to us and ; (iii) Our liability under this guarantee will arise only if we receive a claim in writing , in accordance with the terms of this Guarantee , from on or before 30-05-2023 .'
$d.Paragraphs(4).Range.Text = 'Dated this day of 2022 Signature of the authorized person for and on behalf of the Impactsure Bank Limited PLACE .... DATE ......'
